$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure A1 keeps the same label text (shared string content unchanged, only
# new duplicate shared-string entries get produced as Excel rewrites the table).
$ws.Range("A1").Value = "HK_G_acc_LT"

# New values for A2:A49 (new model / new dataset results).
$values = @(
    86.141304347826093,
    85.326086956521735,
    86.41304347826086,
    78.804347826086953,
    79.891304347826093,
    81.25,
    89.945652173913047,
    88.858695652173907,
    88.58695652173914,
    88.858695652173907,
    73.91304347826086,
    80.706521739130437,
    87.228260869565219,
    86.41304347826086,
    86.956521739130437,
    74.728260869565219,
    76.902173913043484,
    79.619565217391312,
    87.228260869565219,
    88.58695652173914,
    88.315217391304344,
    76.08695652173914,
    77.173913043478265,
    75.815217391304344,
    86.141304347826093,
    83.967391304347828,
    83.152173913043484,
    82.608695652173907,
    79.891304347826093,
    83.152173913043484,
    81.521739130434781,
    83.423913043478265,
    82.065217391304344,
    80.706521739130437,
    81.25,
    73.91304347826086,
    81.521739130434781,
    79.076086956521735,
    79.076086956521735,
    83.967391304347828,
    86.41304347826086,
    84.782608695652172,
    85.054347826086953,
    86.141304347826093,
    85.869565217391312,
    79.076086956521735,
    75,
    82.608695652173907
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
